$wb = $excel.ActiveWorkbook
$wsTotal = $wb.Worksheets.Item(1)   # "总计"
$wsQ = $wb.Worksheets.Item(2)       # currently "2022-Q2", will become "2022-Q3"

# ---------------------------------------------------------------------------
# 1) Before touching the existing quarter sheet, spin off a brand-new sheet
#    that preserves the CURRENT (old Q2) content+formatting verbatim. It is
#    inserted right after the current quarter sheet. (Renaming happens after
#    the source sheet has given up the "2022-Q2" name, to avoid a clash.)
# ---------------------------------------------------------------------------
$wsOldQ2 = $wb.Worksheets.Add($null, $wsQ)

$wsQ.Range("B1:H1").Copy()
$wsOldQ2.Range("B1").PasteSpecial(-4122)
$wsQ.Range("B1:H1").Copy()
$wsOldQ2.Range("B1").PasteSpecial(-4163)

$wsQ.Range("A2:H3").Copy()
$wsOldQ2.Range("A2").PasteSpecial(-4122)
$wsQ.Range("A2:H3").Copy()
$wsOldQ2.Range("A2").PasteSpecial(-4163)

# ---------------------------------------------------------------------------
# 2) Rename the original quarter sheet to "2022-Q3", free up the "2022-Q2"
#    name for the newly split-off sheet, then replace the Q3 sheet's
#    contents with the new quarter's fund data.
# ---------------------------------------------------------------------------
$wsQ.Name = "2022-Q3"
$wsOldQ2.Name = "2022-Q2"
$wsQ.Cells.Clear()

# Re-apply the bold/centered/bordered header style (same style already used
# on the "总计" sheet) to the header row and to the A-column index cells.
$wsTotal.Range("B1").Copy()
$wsQ.Range("B1:H1").PasteSpecial(-4122)
$wsTotal.Range("A2").Copy()
$wsQ.Range("A2:A3").PasteSpecial(-4122)

$wsQ.Range("B1").Value = "基金代码"
$wsQ.Range("C1").Value = "基金名称"
$wsQ.Range("D1").Value = "基金规模"
$wsQ.Range("E1").Value = "股票总仓位"
$wsQ.Range("F1").Value = "仓位占比"
$wsQ.Range("G1").Value = "持有市值(亿元)"
$wsQ.Range("H1").Value = "仓位排名"

$wsQ.Range("A2").Value = 0
$wsQ.Range("B2").Value = "'009619"
$wsQ.Range("C2").Value = "博时女性消费主题混合A"
$wsQ.Range("D2").Value = "'0.54"
$wsQ.Range("E2").Value = "'67.16"
$wsQ.Range("F2").Value = "'3.70"
$wsQ.Range("G2").Value = "'0.0200"
$wsQ.Range("H2").Value = 7

$wsQ.Range("A3").Value = 1
$wsQ.Range("B3").Value = "'009620"
$wsQ.Range("C3").Value = "博时女性消费主题混合C"
$wsQ.Range("D3").Value = "'0.01"
$wsQ.Range("E3").Value = "'67.16"
$wsQ.Range("F3").Value = "'3.70"
$wsQ.Range("G3").Value = "'0.0004"
$wsQ.Range("H3").Value = 7

# ---------------------------------------------------------------------------
# 3) Update the "总计" (totals) sheet: row 2 now reports the Q3 figures and a
#    new row 3 keeps the old Q2 figures that used to live in row 2.
# ---------------------------------------------------------------------------
$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2022-Q2"
$wsTotal.Range("C3").Value = 2
$wsTotal.Range("D3").Value = 0.03

$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("D2").Value = 0.02

$wsTotal.Range("A2").Copy()
$wsTotal.Range("A3").PasteSpecial(-4122)
